$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data per latest scrape.
# Cells holding plain decimal-looking numbers are entered with a
# leading apostrophe so Excel keeps them as text (matching the
# source data format), same as the rest of the untouched column.

$ws.Range("D2").Value = "27.018.87"
$ws.Range("D3").Value = "1.683.88"
$ws.Range("E3").Value = "  +0.84%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'216.24"
$ws.Range("E5").Value = "  +0.25%  "
$ws.Range("D6").Value = "'0.518"
$ws.Range("E6").Value = "  -2.46%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("B8").Value = "Solana"
$ws.Range("C8").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D8").Value = "'21.60"
$ws.Range("E8").Value = "  +6.38%  "
$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").Value = "'0.253"
$ws.Range("E9").Value = "  -0.73%  "
$ws.Range("D10").Value = "'0.0622"
$ws.Range("E10").Value = "  +0.52%  "
$ws.Range("D11").Value = "'0.0889"
$ws.Range("E11").Value = "  -0.54%  "
$ws.Range("D12").Value = "1.923.67"
$ws.Range("E12").Value = "  +0.93%  "
$ws.Range("D13").Value = "1.696.62"
$ws.Range("E13").Value = "  +1.33%  "
$ws.Range("D14").Value = "'4.11"
$ws.Range("E14").Value = "  +0.37%  "
$ws.Range("E15").Value = "  +1.92%  "
$ws.Range("D16").Value = "'66.33"
$ws.Range("E16").Value = "  +0.94%  "
$ws.Range("D17").Value = "'8.22"
$ws.Range("E17").Value = "  +5.44%  "
$ws.Range("D18").Value = "27.049.47"
$ws.Range("E18").Value = "  +0.55%  "
$ws.Range("D19").Value = "'237.08"
$ws.Range("E19").Value = "  +2.01%  "
$ws.Range("D20").Value = "0.0₃0737"
$ws.Range("E20").Value = "  +0.53%  "
$ws.Range("E21").Value = "  +0.06%  "
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("D23").Value = "'9.26"
$ws.Range("E23").Value = "  +0.61%  "
$ws.Range("D24").Value = "'2.12"
$ws.Range("E24").Value = "  -3.95%  "
$ws.Range("D25").Value = "'147.01"
$ws.Range("E25").Value = "  +1.06%  "
$ws.Range("E26").Value = "  +1.33%  "
$ws.Range("D27").Value = "'16.47"
$ws.Range("E27").Value = "  +3.20%  "
$ws.Range("E28").Value = "  -2.87%  "
$ws.Range("E29").Value = "  -0.02%  "
$ws.Range("D30").Value = "'0.0500"
$ws.Range("E30").Value = "  +0.54%  "
$ws.Range("E31").Value = "  -0.04%  "
$ws.Range("E32").Value = "  +0.25%  "
$ws.Range("D33").Value = "1.518.89"
$ws.Range("E33").Value = "  +3.66%  "
$ws.Range("E34").Value = "  +0.71%  "
$ws.Range("D35").Value = "'1.70"
$ws.Range("E35").Value = "  +4.76%  "
$ws.Range("E36").Value = "  -0.39%  "
$ws.Range("D37").Value = "'0.591"
$ws.Range("E37").Value = "  +3.62%  "
$ws.Range("D38").Value = "'0.920"
$ws.Range("E38").Value = "  +2.33%  "
$ws.Range("D39").Value = "'0.0175"
$ws.Range("E39").Value = "  +3.99%  "
$ws.Range("E40").Value = "  +7.53%  "
$ws.Range("D41").Value = "'5.77"
$ws.Range("E41").Value = "  -1.52%  "
$ws.Range("E42").Value = "  -0.02%  "
$ws.Range("D43").Value = "'68.16"
$ws.Range("E43").Value = "  +3.57%  "
$ws.Range("D44").Value = "'2.28"
$ws.Range("E44").Value = "  -0.82%  "
$ws.Range("D45").Value = "1.827.43"
$ws.Range("E46").Value = "  +0.15%  "
$ws.Range("D47").Value = "'90.33"
$ws.Range("E47").Value = "  -0.13%  "
$ws.Range("E48").Value = "  +4.31%  "
$ws.Range("E49").Value = "  -0.08%  "
$ws.Range("D50").Value = "'7.95"
$ws.Range("E50").Value = "  +4.34%  "
$ws.Range("D51").Value = "'0.0507"
$ws.Range("E51").Value = "  -0.32%  "
